# B6-PowerPoint.pptx edit script
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    the custom "Table_0" style ({34EDF4A3-F72F-4920-A81E-442A3133AD4C})
#    to {765252B0-547D-41EC-86C6-5642E5BFF2B3}.
# 2) The presentation theme's colour scheme is repainted from the
#    "Integral" / "Red Violet" palette to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style swap -----------------------------------------------
$newStyleId = "{765252B0-547D-41EC-86C6-5642E5BFF2B3}"
$tableSlides = @(14, 15, 16)
foreach ($idx in $tableSlides) {
    $slide = $p.Slides.Item($idx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour scheme swap ("Integral" -> "Office") --------------
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
